$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Alex Jose Rodriguez Taveras", "Ashlee Ramirez Rosario", "2024-07-11 20:17:06"),
    @("Oscar Daniel Tuletta Mercedes", "Jamil Guzman Feliz", "2024-07-11 20:22:52"),
    @("Oscar Daniel Tuletta Mercedes", "Rafael Antonio Urbaez Hernandez", "2024-07-11 20:32:27"),
    @("Cyd Marie Jorge Chapman", "Edison Yadir Rossis", "2024-07-11 20:33:23"),
    @("Yoelmi Alexander Alcala Valdez", "Yadianna Vargas Pimentel", "2024-07-11 20:33:50")
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $row++
}

$ws.Range("H4").Select()
